$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1143.8889
$ws.Range("I2").Value = 637
$ws.Range("K2").Value = 637
$ws.Range("M2").Value = -524

$ws.Range("H55").Value = 287
$ws.Range("J55").Value = 355
$ws.Range("L55").Value = 355
$ws.Range("N55").Value = -783

$ws.Range("H70").Value = 611500.25
$ws.Range("J70").Value = 1654.2222
$ws.Range("L70").Value = 4962.6666
$ws.Range("N70").Value = -5502.6666

$ws.Range("H73").Value = 611500.25
$ws.Range("J73").Value = 1654.2222
$ws.Range("L73").Value = 4962.6666
$ws.Range("N73").Value = -6834.6666

$ws.Range("H106").Value = 8523.066000000001
$ws.Range("I106").Value = 8953.833000000001
$ws.Range("K106").Value = 8953.833000000001
$ws.Range("M106").Value = -8322.833000000001

$ws.Range("H137").Value = 3171.95
$ws.Range("I137").Value = 2514.1428
$ws.Range("J137").Value = 3526.1538
$ws.Range("K137").Value = 7542.428400000001
$ws.Range("L137").Value = 10578.4614
$ws.Range("M137").Value = -4992.428400000001
$ws.Range("N137").Value = -15678.4614

$ws.Range("H138").Value = 3757.2263
$ws.Range("I138").Value = 1823.4849
$ws.Range("K138").Value = 5470.4547
$ws.Range("M138").Value = -330.4547000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 52635284
$ws.Range("I4").Value = 3806.5334
$ws.Range("K4").Value = 3806.5334
$ws.Range("M4").Value = -3690.5334

$ws.Range("H74").Value = 1585.8
$ws.Range("I74").Value = 1517.5927
$ws.Range("J74").Value = 2199.6667
$ws.Range("K74").Value = 1517.5927
$ws.Range("L74").Value = 2199.6667
$ws.Range("M74").Value = -643.5926999999999
$ws.Range("N74").Value = -3947.6667

$ws.Range("H77").Value = 1585.8
$ws.Range("I77").Value = 1517.5927
$ws.Range("J77").Value = 2199.6667
$ws.Range("K77").Value = 7587.9635
$ws.Range("L77").Value = 10998.3335
$ws.Range("M77").Value = -3219.9635
$ws.Range("N77").Value = -19734.3335

$ws.Range("H88").Value = 2200.2727
$ws.Range("I88").Value = 1809.6666
$ws.Range("J88").Value = 2346.75
$ws.Range("K88").Value = 1809.6666
$ws.Range("L88").Value = 2346.75
$ws.Range("M88").Value = -1403.6666
$ws.Range("N88").Value = -3158.75

$ws.Range("H91").Value = 2200.2727
$ws.Range("I91").Value = 1809.6666
$ws.Range("J91").Value = 2346.75
$ws.Range("K91").Value = 1809.6666
$ws.Range("L91").Value = 2346.75
$ws.Range("M91").Value = -405.6666
$ws.Range("N91").Value = -5154.75

$ws.Range("H122").Value = 2086.3928
$ws.Range("I122").Value = 1824.2
$ws.Range("K122").Value = 5472.6
$ws.Range("M122").Value = -3022.6

$ws.Range("H132").Value = 5558621
$ws.Range("I132").Value = 3046.4
$ws.Range("K132").Value = 9139.200000000001
$ws.Range("M132").Value = -6609.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 412.2857
$ws.Range("J22").Value = 1200
$ws.Range("L22").Value = 1200
$ws.Range("N22").Value = -1900

$ws.Range("H31").Value = 43482256
$ws.Range("I31").Value = 90912370
$ws.Range("K31").Value = 90912370
$ws.Range("M31").Value = -90912075

$ws.Range("H34").Value = 43482256
$ws.Range("I34").Value = 90912370
$ws.Range("K34").Value = 90912370
$ws.Range("M34").Value = -90912168

$ws.Range("H58").Value = 1621.7297
$ws.Range("I58").Value = 1051.9656
$ws.Range("K58").Value = 1051.9656
$ws.Range("M58").Value = -848.9656

$ws.Range("H122").Value = 1971.963
$ws.Range("I122").Value = 1312.2
$ws.Range("K122").Value = 3936.6
$ws.Range("M122").Value = -1486.6

$ws.Range("H132").Value = 3114.4
$ws.Range("I132").Value = 3182.6667
$ws.Range("K132").Value = 9548.000100000001
$ws.Range("M132").Value = -7018.000100000001

$ws.Range("H134").Value = 3605.8
$ws.Range("I134").Value = 3005
$ws.Range("K134").Value = 9015
$ws.Range("M134").Value = -6480

$ws.Range("H136").Value = 1621.7297
$ws.Range("I136").Value = 1051.9656
$ws.Range("K136").Value = 3155.8968
$ws.Range("M136").Value = -605.8968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 520
$ws.Range("I28").Value = 520
$ws.Range("K28").Value = 1560
$ws.Range("M28").Value = -1328

$ws.Range("H45").Value = 509.66666
$ws.Range("J45").Value = 500
$ws.Range("L45").Value = 1500
$ws.Range("N45").Value = -2564

$ws.Range("H52").Value = 430
$ws.Range("J52").Value = 430
$ws.Range("L52").Value = 1290
$ws.Range("N52").Value = -1822

$ws.Range("H62").Value = 14518
$ws.Range("I62").Value = 8299
$ws.Range("K62").Value = 24897
$ws.Range("M62").Value = -24211

$ws.Range("H63").Value = 21643.818
$ws.Range("I63").Value = 8937.5
$ws.Range("K63").Value = 26812.5
$ws.Range("M63").Value = -26063.5

$ws.Range("H65").Value = 14518
$ws.Range("I65").Value = 8299
$ws.Range("K65").Value = 74691
$ws.Range("M65").Value = -71259

$ws.Range("H66").Value = 21643.818
$ws.Range("I66").Value = 8937.5
$ws.Range("K66").Value = 80437.5
$ws.Range("M66").Value = -76693.5

$ws.Range("H87").Value = 99999
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H88").Value = 18230.5
$ws.Range("J88").Value = 18230.5
$ws.Range("L88").Value = 54691.5
$ws.Range("N88").Value = -55547.5

$ws.Range("H90").Value = 99999
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H91").Value = 18230.5
$ws.Range("J91").Value = 18230.5
$ws.Range("L91").Value = 54691.5
$ws.Range("N91").Value = -57655.5

$ws.Range("H113").Value = 826.3333
$ws.Range("I113").Value = 1102
$ws.Range("J113").Value = 688.5
$ws.Range("K113").Value = 3306
$ws.Range("L113").Value = 2065.5
$ws.Range("M113").Value = -1136
$ws.Range("N113").Value = -6405.5

$ws.Range("H123").Value = 7790.2856
$ws.Range("I123").Value = 3533.1667
$ws.Range("K123").Value = 10599.5001
$ws.Range("M123").Value = -8149.500100000001

$ws.Range("H124").Value = 8078.625
$ws.Range("I124").Value = 2259.2
$ws.Range("K124").Value = 6777.599999999999
$ws.Range("M124").Value = -1867.599999999999

$ws.Range("H125").Value = 14604.833
$ws.Range("I125").Value = 9824.75
$ws.Range("J125").Value = 24165
$ws.Range("K125").Value = 29474.25
$ws.Range("L125").Value = 72495
$ws.Range("M125").Value = -24554.25
$ws.Range("N125").Value = -82335

$ws.Range("H126").Value = 12787.667
$ws.Range("I126").Value = 2515
$ws.Range("K126").Value = 7545
$ws.Range("M126").Value = -2605

$ws.Range("H129").Value = 7356555.5
$ws.Range("I129").Value = 13891324
$ws.Range("J129").Value = 4940.75
$ws.Range("K129").Value = 41673972
$ws.Range("L129").Value = 14822.25
$ws.Range("M129").Value = -41668972
$ws.Range("N129").Value = -24822.25

$ws.Range("H139").Value = 5179.9355
$ws.Range("I139").Value = 2681.3333
$ws.Range("K139").Value = 8043.999899999999
$ws.Range("M139").Value = -2903.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4974
$ws.Range("I126").Value = 4662.077
$ws.Range("J126").Value = 5649.8335
$ws.Range("K126").Value = 13986.231
$ws.Range("L126").Value = 16949.5005
$ws.Range("M126").Value = -11516.231
$ws.Range("N126").Value = -21889.5005

$ws.Range("H132").Value = 11116854
$ws.Range("I132").Value = 6781.3335
$ws.Range("J132").Value = 33337000
$ws.Range("K132").Value = 20344.0005
$ws.Range("L132").Value = 100011000
$ws.Range("M132").Value = -17814.0005
$ws.Range("N132").Value = -100016060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5437.9614
$ws.Range("I40").Value = 4539.5
$ws.Range("K40").Value = 4539.5
$ws.Range("M40").Value = -4403.5

$ws.Range("H55").Value = 1554.0555
$ws.Range("I55").Value = 1157.5
$ws.Range("J55").Value = 2049.75
$ws.Range("K55").Value = 1157.5
$ws.Range("L55").Value = 2049.75
$ws.Range("M55").Value = -984.5
$ws.Range("N55").Value = -2395.75

$ws.Range("H133").Value = 74999.60000000001
$ws.Range("J133").Value = 74999.60000000001
$ws.Range("L133").Value = 74999.60000000001
$ws.Range("N133").Value = -80059.60000000001

$ws.Range("H136").Value = 4306.591
$ws.Range("I136").Value = 2375.1428
$ws.Range("J136").Value = 7686.625
$ws.Range("K136").Value = 7125.428400000001
$ws.Range("L136").Value = 23059.875
$ws.Range("M136").Value = -4575.428400000001
$ws.Range("N136").Value = -28159.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 99999
$ws.Range("J95").Value = 99999
$ws.Range("L95").Value = 99999
$ws.Range("N95").Value = -105491

$ws.Range("H126").Value = 6721.875
$ws.Range("I126").Value = 7198.5713
$ws.Range("K126").Value = 21595.7139
$ws.Range("M126").Value = -19125.7139

$ws.Range("H132").Value = 1670040.8
$ws.Range("I132").Value = 3943.3333
$ws.Range("K132").Value = 11829.9999
$ws.Range("M132").Value = -9299.999899999999

$ws.Range("H136").Value = 280200.8
$ws.Range("I136").Value = 2663.4614
$ws.Range("J136").Value = 1001797.9
$ws.Range("K136").Value = 7990.3842
$ws.Range("L136").Value = 3005393.7
$ws.Range("M136").Value = -5440.3842
$ws.Range("N136").Value = -3010493.7
